$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: make room for the new rows -------------------------------
# Old rows 72-76 (5 blank-ish rows) become new rows 72-79 (8 rows), i.e.
# we need 3 extra rows. Insert them right after the existing row 76 so the
# existing row 72-76 formatting stays put (and gets overwritten below) and
# rows 77-79 appear as brand-new rows, matching the diff structure.
$ws.Rows("77:79").Insert()

# --- Step 2: set cell VALUES in the exact chronological order the
# original author must have typed them in, so new shared-string table
# entries land on the same indices as the target file (103..108).
$ws.Range("A73").Value2 = "closets_page.dart"
$ws.Range("B74").Value2 = "Tạo tủ đồ mới"
$ws.Range("B71").Value2 = "Tủ đồ của bạn đang trống. Hãy thêm đồ vật trước nhé!"
$ws.Range("A74").Value2 = "Create new closet"
$ws.Range("A75").Value2 = "Closet name"
$ws.Range("B75").Value2 = "Tên tủ đồ"
$ws.Range("A76").Value2 = "Cancel"
$ws.Range("B76").Value2 = "Hủy"
$ws.Range("A77").Value2 = "Save"
$ws.Range("B77").Value2 = "Lưu"

# --- Step 3: formatting -------------------------------------------------
# Row 72 becomes an empty spacer row with no border (new style).
$ws.Range("A71:B71").Copy()
$ws.Range("A72:B72").PasteSpecial(-4122)
$ws.Range("A72:B72").Borders.LineStyle = 0

# Row 73 is a new section header ("closets_page.dart") -> bold/yellow style.
$ws.Range("A69:B69").Copy()
$ws.Range("A73:B73").PasteSpecial(-4122)

# Row 74 ("Create new closet" / "Tạo tủ đồ mới") -> plain bordered style.
$ws.Range("A2:B2").Copy()
$ws.Range("A74:B74").PasteSpecial(-4122)

# Row 75 ("Closet name" / "Tên tủ đồ") -> bordered style variant.
$ws.Range("A70:B70").Copy()
$ws.Range("A75:B75").PasteSpecial(-4122)

# Row 76 ("Cancel" / "Hủy") -> plain bordered style.
$ws.Range("A2:B2").Copy()
$ws.Range("A76:B76").PasteSpecial(-4122)

# Row 77 ("Save" / "Lưu") -> plain bordered style.
$ws.Range("A28:B28").Copy()
$ws.Range("A77:B77").PasteSpecial(-4122)

# Row 78 (blank) -> same style family as row 77.
$ws.Range("A28:B28").Copy()
$ws.Range("A78:B78").PasteSpecial(-4122)

# Row 79 (blank) -> bordered style variant (same family as row 75).
$ws.Range("A70:B70").Copy()
$ws.Range("A79:B79").PasteSpecial(-4122)

# --- Step 4: merge the new header cell, like the other section headers --
$ws.Range("A73:B73").Merge()

# --- Step 5: tidy up selection / view, matching the final workbook state
$ws.Range("A78").Select()
$ws.Application.ActiveWindow.ScrollRow = 65
